$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3734442.8
$ws.Range("I132").Value = 3005.6897
$ws.Range("J132").Value = 27781482
$ws.Range("K132").Value = 9017.069100000001
$ws.Range("L132").Value = 83344446
$ws.Range("M132").Value = -6487.069100000001
$ws.Range("N132").Value = -83349506

$ws.Range("H137").Value = 6897387
$ws.Range("I137").Value = 594.3684
$ws.Range("J137").Value = 20001292
$ws.Range("K137").Value = 1783.1052
$ws.Range("L137").Value = 60003876
$ws.Range("M137").Value = 766.8948
$ws.Range("N137").Value = -60008976

$ws.Range("H138").Value = 6668195.5
$ws.Range("I138").Value = 9804858
$ws.Range("J138").Value = 2787.5
$ws.Range("K138").Value = 29414574
$ws.Range("L138").Value = 8362.5
$ws.Range("M138").Value = -29409434
$ws.Range("N138").Value = -18642.5

$ws.Range("H141").Value = 1035.1277
$ws.Range("I141").Value = 983.8049
$ws.Range("J141").Value = 1385.8334
$ws.Range("K141").Value = 2951.4147
$ws.Range("L141").Value = 4157.5002
$ws.Range("M141").Value = 2228.5853
$ws.Range("N141").Value = -14517.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9189.651
$ws.Range("I32").Value = 9238.367
$ws.Range("K32").Value = 9238.367
$ws.Range("M32").Value = -8951.367

$ws.Range("H61").Value = 9805317
$ws.Range("I61").Value = 10418060
$ws.Range("J61").Value = 1433.6666
$ws.Range("K61").Value = 10418060
$ws.Range("L61").Value = 1433.6666
$ws.Range("M61").Value = -10417848
$ws.Range("N61").Value = -1857.6666

$ws.Range("H74").Value = 7464400
$ws.Range("I74").Value = 8621858
$ws.Range("J74").Value = 5223.778
$ws.Range("K74").Value = 8621858
$ws.Range("L74").Value = 5223.778
$ws.Range("M74").Value = -8620984
$ws.Range("N74").Value = -6971.778

$ws.Range("H77").Value = 7464400
$ws.Range("I77").Value = 8621858
$ws.Range("J77").Value = 5223.778
$ws.Range("K77").Value = 43109290
$ws.Range("L77").Value = 26118.89
$ws.Range("M77").Value = -43104922
$ws.Range("N77").Value = -34854.89

$ws.Range("H122").Value = 5685.3105
$ws.Range("I122").Value = 7549.2104
$ws.Range("J122").Value = 2143.9
$ws.Range("K122").Value = 22647.6312
$ws.Range("L122").Value = 6431.700000000001
$ws.Range("M122").Value = -20197.6312
$ws.Range("N122").Value = -11331.7

$ws.Range("H132").Value = 3473776
$ws.Range("I132").Value = 4465589
$ws.Range("J132").Value = 2429.625
$ws.Range("K132").Value = 13396767
$ws.Range("L132").Value = 7288.875
$ws.Range("M132").Value = -13394237
$ws.Range("N132").Value = -12348.875

$ws.Range("H136").Value = 9805317
$ws.Range("I136").Value = 10418060
$ws.Range("J136").Value = 1433.6666
$ws.Range("K136").Value = 31254180
$ws.Range("L136").Value = 4300.9998
$ws.Range("M136").Value = -31251630
$ws.Range("N136").Value = -9400.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 11350.889
$ws.Range("I80").Value = 40214
$ws.Range("J80").Value = 249.6923
$ws.Range("K80").Value = 40214
$ws.Range("L80").Value = 249.6923
$ws.Range("M80").Value = -39216
$ws.Range("N80").Value = -2245.6923

$ws.Range("H83").Value = 11350.889
$ws.Range("I83").Value = 40214
$ws.Range("J83").Value = 249.6923
$ws.Range("K83").Value = 201070
$ws.Range("L83").Value = 1248.4615
$ws.Range("M83").Value = -196078
$ws.Range("N83").Value = -11232.4615

$ws.Range("H134").Value = 2463.2593
$ws.Range("I134").Value = 1333.1842
$ws.Range("J134").Value = 5147.1875
$ws.Range("K134").Value = 3999.5526
$ws.Range("L134").Value = 15441.5625
$ws.Range("M134").Value = -1464.5526
$ws.Range("N134").Value = -20511.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1042.3265
$ws.Range("I58").Value = 501.82758
$ws.Range("J58").Value = 1826.05
$ws.Range("K58").Value = 501.82758
$ws.Range("L58").Value = 1826.05
$ws.Range("M58").Value = -298.82758
$ws.Range("N58").Value = -2232.05

$ws.Range("H132").Value = 11113128
$ws.Range("I132").Value = 14287470
$ws.Range("J132").Value = 2934
$ws.Range("K132").Value = 42862410
$ws.Range("L132").Value = 8802
$ws.Range("M132").Value = -42859880
$ws.Range("N132").Value = -13862

$ws.Range("H134").Value = 426494.3
$ws.Range("I134").Value = 1312.25
$ws.Range("J134").Value = 1985495.1
$ws.Range("K134").Value = 3936.75
$ws.Range("L134").Value = 5956485.300000001
$ws.Range("M134").Value = -1401.75
$ws.Range("N134").Value = -5961555.300000001

$ws.Range("H136").Value = 1042.3265
$ws.Range("I136").Value = 501.82758
$ws.Range("J136").Value = 1826.05
$ws.Range("K136").Value = 1505.48274
$ws.Range("L136").Value = 5478.15
$ws.Range("M136").Value = 1044.51726
$ws.Range("N136").Value = -10578.15

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 915.8333
$ws.Range("I86").Value = 323.75
$ws.Range("J86").Value = 2100
$ws.Range("K86").Value = 971.25
$ws.Range("L86").Value = 6300
$ws.Range("M86").Value = 214.75
$ws.Range("N86").Value = -8672

$ws.Range("H89").Value = 915.8333
$ws.Range("I89").Value = 323.75
$ws.Range("J89").Value = 2100
$ws.Range("K89").Value = 2913.75
$ws.Range("L89").Value = 18900
$ws.Range("M89").Value = 3014.25
$ws.Range("N89").Value = -30756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4140.409
$ws.Range("I102").Value = 4199.7617
$ws.Range("K102").Value = 4199.7617
$ws.Range("M102").Value = -2577.7617

$ws.Range("H122").Value = 2471026.5
$ws.Range("J122").Value = 2265.818
$ws.Range("L122").Value = 6797.454000000001
$ws.Range("N122").Value = -11697.454

$ws.Range("H132").Value = 3218.3062
$ws.Range("I132").Value = 2341.6287
$ws.Range("J132").Value = 5410
$ws.Range("K132").Value = 7024.886100000001
$ws.Range("L132").Value = 16230
$ws.Range("M132").Value = -4494.886100000001
$ws.Range("N132").Value = -21290

$ws.Range("H136").Value = 19186.285
$ws.Range("J136").Value = 19186.285
$ws.Range("L136").Value = 57558.855
$ws.Range("N136").Value = -62658.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6306.5454
$ws.Range("I7").Value = 7366.143
$ws.Range("J7").Value = 5525.7896
$ws.Range("K7").Value = 7366.143
$ws.Range("L7").Value = 5525.7896
$ws.Range("M7").Value = -7254.143
$ws.Range("N7").Value = -5749.7896

$ws.Range("H40").Value = 5248.4287
$ws.Range("I40").Value = 7700
$ws.Range("K40").Value = 7700
$ws.Range("M40").Value = -7564

$ws.Range("H61").Value = 1538.4615
$ws.Range("I61").Value = 1356.4286
$ws.Range("J61").Value = 1750.8334
$ws.Range("K61").Value = 1356.4286
$ws.Range("L61").Value = 1750.8334
$ws.Range("M61").Value = -1154.4286
$ws.Range("N61").Value = -2154.8334

$ws.Range("H82").Value = 2430.3
$ws.Range("J82").Value = 2400.4285
$ws.Range("L82").Value = 2400.4285
$ws.Range("N82").Value = -3122.4285

$ws.Range("H85").Value = 2430.3
$ws.Range("J85").Value = 2400.4285
$ws.Range("L85").Value = 2400.4285
$ws.Range("N85").Value = -4896.4285

$ws.Range("H113").Value = 1538.4615
$ws.Range("I113").Value = 1356.4286
$ws.Range("J113").Value = 1750.8334
$ws.Range("K113").Value = 1356.4286
$ws.Range("L113").Value = 1750.8334
$ws.Range("M113").Value = 813.5714
$ws.Range("N113").Value = -6090.8334

$ws.Range("H122").Value = 4579.439
$ws.Range("I122").Value = 4309.885
$ws.Range("J122").Value = 5046.6665
$ws.Range("K122").Value = 12929.655
$ws.Range("L122").Value = 15139.9995
$ws.Range("M122").Value = -10479.655
$ws.Range("N122").Value = -20039.9995

$ws.Range("H126").Value = 6306.5454
$ws.Range("I126").Value = 7366.143
$ws.Range("J126").Value = 5525.7896
$ws.Range("K126").Value = 22098.429
$ws.Range("L126").Value = 16577.3688
$ws.Range("M126").Value = -19628.429
$ws.Range("N126").Value = -21517.3688

$ws.Range("H132").Value = 6029398.5
$ws.Range("I132").Value = 2988.4531
$ws.Range("J132").Value = 26328886
$ws.Range("K132").Value = 8965.3593
$ws.Range("L132").Value = 78986658
$ws.Range("M132").Value = -6435.3593
$ws.Range("N132").Value = -78991718

$ws.Range("H136").Value = 8477385
$ws.Range("I136").Value = 11112081
$ws.Range("J136").Value = 8717.5
$ws.Range("K136").Value = 33336243
$ws.Range("L136").Value = 26152.5
$ws.Range("M136").Value = -33333693
$ws.Range("N136").Value = -31252.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12780.8
$ws.Range("I62").Value = 6200.5
$ws.Range("J62").Value = 17167.666
$ws.Range("K62").Value = 6200.5
$ws.Range("L62").Value = 17167.666
$ws.Range("M62").Value = -5576.5
$ws.Range("N62").Value = -18415.666

$ws.Range("H65").Value = 12780.8
$ws.Range("I65").Value = 6200.5
$ws.Range("J65").Value = 17167.666
$ws.Range("K65").Value = 31002.5
$ws.Range("L65").Value = 85838.33
$ws.Range("M65").Value = -27882.5
$ws.Range("N65").Value = -92078.33

$ws.Range("H122").Value = 1989.1111
$ws.Range("I122").Value = 2345.611
$ws.Range("J122").Value = 1276.1111
$ws.Range("K122").Value = 7036.833
$ws.Range("L122").Value = 3828.3333
$ws.Range("M122").Value = -4586.833
$ws.Range("N122").Value = -8728.3333

$ws.Range("H136").Value = 800.2222
$ws.Range("I136").Value = 655.2826
$ws.Range("J136").Value = 1633.625
$ws.Range("K136").Value = 1965.8478
$ws.Range("L136").Value = 4900.875
$ws.Range("M136").Value = 584.1522
$ws.Range("N136").Value = -10000.875
